# ---------------------------------------------------------------------
# Stand-Up.docx edit:
#   1) Split the run containing "sqflite" in the 11-11-2021 Summary
#      paragraph, wrapping the word in proofErr spellStart/spellEnd
#      (matches Word's "flag as a spelling suggestion" markup).
#   2) Do the same for "WinHistoryCard" in the 12-11-2021 Summary
#      paragraph, and append a new 13-11-2021 Date/TODO/Summary block
#      right after it.
#
# InsertXML on a Range *replaces that range's contents*, so each edit
# below selects the *entire* paragraph (Paragraph.Range, which includes
# the trailing paragraph mark) and supplies the full replacement
# paragraph(s) as WordprocessingML, wrapped in the mandatory
# pkg:package envelope.
# ---------------------------------------------------------------------

$d = $word.ActiveDocument

function New-PkgXml([string]$bodyXml) {
    return @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
$bodyXml
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
}

# Locate the two paragraphs by their (unique) current text, rather than
# hard-coded indices, so the script is resilient to minor surprises.
$sqfliteRange = $d.Content
$sqfliteRange.Find.ClearFormatting()
$sqfliteRange.Find.Execute("sqflite", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$paraSqflite = $sqfliteRange.Paragraphs(1)

$winCardRange = $d.Content
$winCardRange.Find.ClearFormatting()
$winCardRange.Find.Execute("WinHistoryCard", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$paraWinCard = $winCardRange.Paragraphs(1)

# -----------------------------------------------------------------
# 1) "...backend and sqflite to make history screen..."
# -----------------------------------------------------------------
$para1Xml = @"
<w:p w14:paraId="5168750A" w14:textId="753F0DA2" w:rsidR="00897E17" w:rsidRDefault="00897E17" w:rsidP="00E225FD">
  <w:pPr>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve">Summary: </w:t>
  </w:r>
  <w:r w:rsidR="00A9277E">
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve">Today I&#8217;ve been working on the design </w:t>
  </w:r>
  <w:r w:rsidR="00FE39E0">
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>section,</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve"> and I&#8217;ve finished it for now, and I worked on backend and </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>sqflite</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve"> to make history screen</w:t>
  </w:r>
  <w:r w:rsidR="000C1CFE">
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>, and it took 3:45 hours from me to complete it</w:t>
  </w:r>
  <w:r w:rsidR="00A9277E">
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>.</w:t>
  </w:r>
</w:p>
"@
$paraSqflite.Range.InsertXML((New-PkgXml $para1Xml))

# -----------------------------------------------------------------
# 2) "...completed the WinHistoryCard and that will show cards..."
#    plus the brand-new 13-11-2021 Date / TODO / Summary paragraphs.
# -----------------------------------------------------------------
$para2Xml = @"
<w:p w14:paraId="008F84A3" w14:textId="3995A928" w:rsidR="00425089" w:rsidRPr="00285C95" w:rsidRDefault="00425089" w:rsidP="00E225FD">
  <w:pPr>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:lang w:val="en-US" w:bidi="ar-SY"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:lang w:val="en-US" w:bidi="ar-SY"/>
    </w:rPr>
    <w:t xml:space="preserve">Summary: </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:lang w:val="en-US" w:bidi="ar-SY"/>
    </w:rPr>
    <w:t xml:space="preserve">Today I&#8217;ve completed the </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:lang w:val="en-US" w:bidi="ar-SY"/>
    </w:rPr>
    <w:t>WinHistoryCard</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:lang w:val="en-US" w:bidi="ar-SY"/>
    </w:rPr>
    <w:t xml:space="preserve"> and that will show cards as a history instead of text.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve">Date: </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>13-11-2021 &#8211; 5 hours.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve">TODO: </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>Complete the project and write documentation.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve">Summary: </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>Today I&#8217;ve finished the project and wrote the documentation.</w:t>
  </w:r>
</w:p>
"@
$paraWinCard.Range.InsertXML((New-PkgXml $para2Xml))
